$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: title/link update
$ws.Range("D10").Value = "Terraform Cloud 권한 관리하기 (feat. Vault, AWS)"
$ws.Range("E10").Value = "https://rokrokss.com/post/2020/04/14/Terraform-Cloud-%EA%B6%8C%ED%95%9C-%EA%B4%80%EB%A6%AC%ED%95%98%EA%B8%B0-Vault-AWS.html"

# Row 26: title update
$ws.Range("D26").Value = "생성 모델의 새로운 흐름 확산 모델(Diffusion model)에 관하여"

# Row 51: title/link update
$ws.Range("D51").Value = "[css] input 박스 또는 textarea 박스 클릭했을 때 테두리 강조되지 않게 하기(:focus 의사 클래스)"
$ws.Range("E51").Value = "https://bskyvision.com/1256"
